$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.03
$ws.Range("C2").Value = 0.31

$ws.Range("B3").Value = 89.11

$ws.Range("B4").Value = 169.11

$ws.Range("B5").Value = 0.24
$ws.Range("C5").Value = 0.62

$ws.Range("B6").Value = 20.54

$ws.Range("B7").Value = 1.19
$ws.Range("C7").Value = 0.28

$ws.Range("B8").Value = 62.03
